$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 6323
$ws.Range("I40").Value = 4462.4165
$ws.Range("J40").Value = 7636.353
$ws.Range("K40").Value = 4462.4165
$ws.Range("L40").Value = 7636.353
$ws.Range("M40").Value = -4287.4165
$ws.Range("N40").Value = -7986.353
$ws.Range("H98").Value = 1009.55554
$ws.Range("I98").Value = 1009.55554
$ws.Range("K98").Value = 1009.55554
$ws.Range("M98").Value = 488.44446
$ws.Range("H100").Value = 4158.706
$ws.Range("I100").Value = 2375.111
$ws.Range("K100").Value = 2375.111
$ws.Range("M100").Value = -1834.111
$ws.Range("H113").Value = 4899.4707
$ws.Range("I113").Value = 3180.6
$ws.Range("K113").Value = 3180.6
$ws.Range("M113").Value = 73.40000000000009
$ws.Range("H122").Value = 1009.55554
$ws.Range("I122").Value = 1009.55554
$ws.Range("K122").Value = 3028.66662
$ws.Range("M122").Value = -578.66662

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H43").Value = 6290102
$ws.Range("J43").Value = 483503
$ws.Range("L43").Value = 483503
$ws.Range("N43").Value = -484129
$ws.Range("H61").Value = 4729.1177
$ws.Range("I61").Value = 3101.6
$ws.Range("K61").Value = 3101.6
$ws.Range("M61").Value = -2889.6
$ws.Range("H96").Value = 22500
$ws.Range("J96").Value = 20000
$ws.Range("L96").Value = 20000
$ws.Range("N96").Value = -25492
$ws.Range("H136").Value = 4729.1177
$ws.Range("I136").Value = 3101.6
$ws.Range("K136").Value = 9304.799999999999
$ws.Range("M136").Value = -6754.799999999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 1580
$ws.Range("I105").Value = 1580
$ws.Range("K105").Value = 1580
$ws.Range("M105").Value = 167
$ws.Range("H134").Value = 2364.05
$ws.Range("I134").Value = 2248.9443
$ws.Range("J134").Value = 3400
$ws.Range("K134").Value = 6746.8329
$ws.Range("L134").Value = 10200
$ws.Range("M134").Value = -4211.8329
$ws.Range("N134").Value = -15270
$ws.Range("H135").Value = 197497.5
$ws.Range("J135").Value = 197497.5
$ws.Range("L135").Value = 197497.5
$ws.Range("N135").Value = -207637.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 200
$ws.Range("J22").Value = 200
$ws.Range("L22").Value = 200
$ws.Range("N22").Value = -900
$ws.Range("H31").Value = 4384.855
$ws.Range("I31").Value = 2905.3125
$ws.Range("K31").Value = 2905.3125
$ws.Range("M31").Value = -2610.3125
$ws.Range("H34").Value = 4384.855
$ws.Range("I34").Value = 2905.3125
$ws.Range("K34").Value = 2905.3125
$ws.Range("M34").Value = -2703.3125
$ws.Range("H41").Value = 32599.814
$ws.Range("I41").Value = 3950
$ws.Range("J41").Value = 34891.8
$ws.Range("K41").Value = 3950
$ws.Range("L41").Value = 34891.8
$ws.Range("M41").Value = -3522
$ws.Range("N41").Value = -35747.8
$ws.Range("H58").Value = 2340
$ws.Range("I58").Value = 1827.8182
$ws.Range("J58").Value = 3748.5
$ws.Range("K58").Value = 1827.8182
$ws.Range("L58").Value = 3748.5
$ws.Range("M58").Value = -1624.8182
$ws.Range("N58").Value = -4154.5
$ws.Range("H59").Value = 123551140
$ws.Range("J59").Value = 123551140
$ws.Range("L59").Value = 123551140
$ws.Range("N59").Value = -123553430
$ws.Range("H69").Value = 39999.5
$ws.Range("I69").Value = 39999
$ws.Range("J69").Value = 40000
$ws.Range("K69").Value = 39999
$ws.Range("L69").Value = 40000
$ws.Range("M69").Value = -39250
$ws.Range("N69").Value = -41498
$ws.Range("H72").Value = 39999.5
$ws.Range("I72").Value = 39999
$ws.Range("J72").Value = 40000
$ws.Range("K72").Value = 119997
$ws.Range("L72").Value = 120000
$ws.Range("M72").Value = -116253
$ws.Range("N72").Value = -127488
$ws.Range("H92").Value = 39000
$ws.Range("J92").Value = 39000
$ws.Range("L92").Value = 39000
$ws.Range("N92").Value = -43992
$ws.Range("H99").Value = 2223.5
$ws.Range("I99").Value = 2165.6667
$ws.Range("K99").Value = 2165.6667
$ws.Range("M99").Value = -667.6667000000002
$ws.Range("H126").Value = 2223.5
$ws.Range("I126").Value = 2165.6667
$ws.Range("K126").Value = 6497.000100000001
$ws.Range("M126").Value = -4027.000100000001
$ws.Range("H136").Value = 2340
$ws.Range("I136").Value = 1827.8182
$ws.Range("J136").Value = 3748.5
$ws.Range("K136").Value = 5483.4546
$ws.Range("L136").Value = 11245.5
$ws.Range("M136").Value = -2933.4546
$ws.Range("N136").Value = -16345.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H38").Value = 399
$ws.Range("J38").Value = 0
$ws.Range("L38").Value = 0
$ws.Range("N38").ClearContents()
$ws.Range("H51").Value = 1630.25
$ws.Range("I51").Value = 1630.25
$ws.Range("J51").Value = 0
$ws.Range("K51").Value = 4890.75
$ws.Range("L51").Value = 0
$ws.Range("M51").Value = -4430.75
$ws.Range("N51").ClearContents()
$ws.Range("H80").Value = 3957.25
$ws.Range("J80").Value = 4475.5
$ws.Range("L80").Value = 13426.5
$ws.Range("N80").Value = -15298.5
$ws.Range("H83").Value = 3957.25
$ws.Range("J83").Value = 4475.5
$ws.Range("L83").Value = 40279.5
$ws.Range("N83").Value = -49639.5
$ws.Range("H107").Value = 612.75
$ws.Range("I107").Value = 685.3333
$ws.Range("K107").Value = 2055.9999
$ws.Range("M107").Value = -135.9998999999998
$ws.Range("H128").Value = 256000
$ws.Range("I128").Value = 256000
$ws.Range("K128").Value = 768000
$ws.Range("M128").Value = -763020
$ws.Range("H129").Value = 1881.6111
$ws.Range("J129").Value = 2248.0833
$ws.Range("L129").Value = 6744.249899999999
$ws.Range("N129").Value = -16744.2499

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H11").Value = 10780913
$ws.Range("I11").Value = 11381010
$ws.Range("J11").Value = 8294795.5
$ws.Range("K11").Value = 11381010
$ws.Range("L11").Value = 8294795.5
$ws.Range("M11").Value = -11380871
$ws.Range("N11").Value = -8295073.5
$ws.Range("H20").Value = 0
$ws.Range("J20").Value = 0
$ws.Range("L20").Value = 0
$ws.Range("N20").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 8658
$ws.Range("I7").Value = 8579
$ws.Range("K7").Value = 8579
$ws.Range("M7").Value = -8467
$ws.Range("H20").Value = 15000
$ws.Range("I20").Value = 15000
$ws.Range("J20").Value = 0
$ws.Range("K20").Value = 15000
$ws.Range("L20").Value = 0
$ws.Range("M20").Value = -14774
$ws.Range("N20").ClearContents()
$ws.Range("H22").Value = 3582.3333
$ws.Range("I22").Value = 3112.5
$ws.Range("J22").Value = 3958.2
$ws.Range("K22").Value = 3112.5
$ws.Range("L22").Value = 3958.2
$ws.Range("M22").Value = -2817.5
$ws.Range("N22").Value = -4548.2
$ws.Range("H27").Value = 3582.3333
$ws.Range("I27").Value = 3112.5
$ws.Range("J27").Value = 3958.2
$ws.Range("K27").Value = 3112.5
$ws.Range("L27").Value = 3958.2
$ws.Range("M27").Value = -3005.5
$ws.Range("N27").Value = -4172.2
$ws.Range("H40").Value = 3998.1667
$ws.Range("I40").Value = 3797.8
$ws.Range("K40").Value = 3797.8
$ws.Range("M40").Value = -3661.8
$ws.Range("H100").Value = 9171
$ws.Range("I100").Value = 7789.3335
$ws.Range("K100").Value = 7789.3335
$ws.Range("M100").Value = -7248.3335
$ws.Range("H126").Value = 8658
$ws.Range("I126").Value = 8579
$ws.Range("K126").Value = 25737
$ws.Range("M126").Value = -23267

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H4").Value = 29306.643
$ws.Range("I4").Value = 31557.154
$ws.Range("K4").Value = 31557.154
$ws.Range("M4").Value = -31444.154
$ws.Range("H5").Value = 16000200
$ws.Range("I5").Value = 15000250
$ws.Range("K5").Value = 15000250
$ws.Range("M5").Value = -15000138
$ws.Range("H126").Value = 5750.095
$ws.Range("I126").Value = 3667.75
$ws.Range("K126").Value = 11003.25
$ws.Range("M126").Value = -8533.25
